# garment_data.xlsx — convert the Actual/Target/Variance KPI block from
# mixed percent/thousands figures to plain whole numbers, and tidy up the
# sheet (number format, column widths, selection, page setup) to match how
# Excel left the file after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- new KPI figures (row 2: Plan Vs. Actual, row 3: Efficiency, row 4: Lost Time) ---
$ws.Range("B2").Value = 80
$ws.Range("C2").Value = 100
# D2 already holds =B2-C2 and recalculates on its own

$ws.Range("B3").Value = 65
$ws.Range("C3").Value = 70
# D3 already holds the shared formula =B3-C3

$ws.Range("B4").Value = 120
$ws.Range("C4").Value = 100
# D4 already holds the shared formula =B3-C3 (si="0")

# --- all of the Actual/Target/Variance numbers now share a plain "0" format ---
$ws.Range("B2:D4").NumberFormat = "0"

# --- widen B:D to fit their new contents (mirrors a manual column auto-fit) ---
$ws.Columns("B:D").AutoFit()

# --- page setup tweaks left over from the edit (A4, portrait) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- cursor ends up on D10 ---
$ws.Range("D10").Select() | Out-Null
